# "A long overdue commit of many random changes."
# - Header H3: "Length (inches)" -> "Cut to Length (inches)"
# - H8 corrected from 60.26 to 60.28
# - Widen columns G:H
# - Append new parts (rows 21-31) to the materials list
# - Update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Corrected length for line 5 (Lower long)
$ws.Range("H8").Value = 60.28

# Column G and H both widened to match (target stored width 22.42578125 chars)
$ws.Columns.Item(7).ColumnWidth = 21.59
$ws.Columns.Item(8).ColumnWidth = 21.59

# New hardware / fastener rows appended below the panel section
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 3383
$ws.Range("E21").Value = 64
$ws.Range("G21").Value = "Single tab end fastener"

$ws.Range("C22").Value = 13
$ws.Range("D22").Value = 3098
$ws.Range("E22").Value = 24
$ws.Range("G22").Value = "Double anchor fastener"

$ws.Range("C23").Value = 14
$ws.Range("D23").Value = 2059
$ws.Range("E23").Value = 8
$ws.Range("G23").Value = "Door hanger"

$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 2067
$ws.Range("E24").Value = 8
$ws.Range("G24").Value = "Door glide"

$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 2103
$ws.Range("E25").Value = 16
$ws.Range("G25").Value = "Inter-series hinge"

$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 2066
$ws.Range("E26").Value = 16
$ws.Range("G26").Value = "Within-series hinge"

$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 2079
$ws.Range("E27").Value = 24
$ws.Range("G27").Value = "Offset door handle"

$ws.Range("C28").Value = 19
$ws.Range("D28").Value = 4515
$ws.Range("E28").Value = 16
$ws.Range("G28").Value = "Interseries joining plate"

$ws.Range("C29").Value = 20
$ws.Range("D29").Value = 3321
$ws.Range("E29").Value = 32
$ws.Range("G29").Value = "10 series nut"

$ws.Range("C30").Value = 21
$ws.Range("D30").Value = 3320
$ws.Range("E30").Value = 32
$ws.Range("G30").Value = "15 series nut"

$ws.Range("C31").Value = 22
$ws.Range("D31").Value = 2116
$ws.Range("E31").Value = "170 feet"
$ws.Range("G31").Value = "Rubber panel gasket"

# Header text tweak (added last so the new shared string lands at the end,
# matching the order the workbook's string table was built in)
$ws.Range("H3").Value = "Cut to Length (inches)"

# Leave the selection where the author left it
$ws.Range("I18").Select()
